$wb = $excel.ActiveWorkbook

# Scheduled market-price refresh: update currentAveragePrice / LevePrice / LeveProfit
# columns (H-N) across all job sheets, per the latest price snapshot.

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 55555736
$ws.Range("I33").Value = 111.13333
$ws.Range("J33").Value = 333333860
$ws.Range("K33").Value = 111.13333
$ws.Range("L33").Value = 333333860
$ws.Range("M33").Value = 117.86667
$ws.Range("N33").Value = -333334318

$ws.Range("H40").Value = 3040.4
$ws.Range("I40").Value = 3093.8
$ws.Range("K40").Value = 3093.8
$ws.Range("M40").Value = -2918.8

$ws.Range("H48").Value = 243.35
$ws.Range("I48").Value = 117
$ws.Range("K48").Value = 351
$ws.Range("M48").Value = -59

$ws.Range("H56").Value = 243.35
$ws.Range("I56").Value = 117
$ws.Range("K56").Value = 351
$ws.Range("M56").Value = 183

$ws.Range("H112").Value = 1595.7333
$ws.Range("J112").Value = 1793.4
$ws.Range("L112").Value = 5380.200000000001
$ws.Range("N112").Value = -7596.200000000001

$ws.Range("H135").Value = 1334.6562
$ws.Range("I135").Value = 1300.6897
$ws.Range("K135").Value = 11706.2073
$ws.Range("M135").Value = -9171.207299999998

$ws.Range("H138").Value = 1740.0952
$ws.Range("I138").Value = 660.86957
$ws.Range("J138").Value = 3046.5264
$ws.Range("K138").Value = 1982.60871
$ws.Range("L138").Value = 9139.5792
$ws.Range("M138").Value = 3157.39129
$ws.Range("N138").Value = -19419.5792


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2415.3684
$ws.Range("I63").Value = 3197.0908
$ws.Range("J63").Value = 1340.5
$ws.Range("K63").Value = 3197.0908
$ws.Range("L63").Value = 1340.5
$ws.Range("M63").Value = -2511.0908
$ws.Range("N63").Value = -2712.5

$ws.Range("H66").Value = 2415.3684
$ws.Range("I66").Value = 3197.0908
$ws.Range("J66").Value = 1340.5
$ws.Range("K66").Value = 15985.454
$ws.Range("L66").Value = 6702.5
$ws.Range("M66").Value = -12553.454
$ws.Range("N66").Value = -13566.5

$ws.Range("H74").Value = 3224.739
$ws.Range("I74").Value = 2814.7222
$ws.Range("J74").Value = 4700.8
$ws.Range("K74").Value = 2814.7222
$ws.Range("L74").Value = 4700.8
$ws.Range("M74").Value = -1940.7222
$ws.Range("N74").Value = -6448.8

$ws.Range("H77").Value = 3224.739
$ws.Range("I77").Value = 2814.7222
$ws.Range("J77").Value = 4700.8
$ws.Range("K77").Value = 14073.611
$ws.Range("L77").Value = 23504
$ws.Range("M77").Value = -9705.611000000001
$ws.Range("N77").Value = -32240

$ws.Range("H88").Value = 2862.1177
$ws.Range("I88").Value = 2699.4
$ws.Range("J88").Value = 2929.9167
$ws.Range("K88").Value = 2699.4
$ws.Range("L88").Value = 2929.9167
$ws.Range("M88").Value = -2293.4
$ws.Range("N88").Value = -3741.9167

$ws.Range("H91").Value = 2862.1177
$ws.Range("I91").Value = 2699.4
$ws.Range("J91").Value = 2929.9167
$ws.Range("K91").Value = 2699.4
$ws.Range("L91").Value = 2929.9167
$ws.Range("M91").Value = -1295.4
$ws.Range("N91").Value = -5737.9167

$ws.Range("H122").Value = 7590.4614
$ws.Range("I122").Value = 8246.6
$ws.Range("J122").Value = 5403.3335
$ws.Range("K122").Value = 24739.8
$ws.Range("L122").Value = 16210.0005
$ws.Range("M122").Value = -22289.8
$ws.Range("N122").Value = -21110.0005


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 2506
$ws.Range("I22").Value = 359
$ws.Range("K22").Value = 359
$ws.Range("M22").Value = -186

$ws.Range("H86").Value = 1982.4615
$ws.Range("I86").Value = 1705.9375
$ws.Range("K86").Value = 1705.9375
$ws.Range("M86").Value = -582.9375

$ws.Range("H89").Value = 1982.4615
$ws.Range("I89").Value = 1705.9375
$ws.Range("K89").Value = 8529.6875
$ws.Range("M89").Value = -2913.6875

$ws.Range("H94").Value = 1912.6207
$ws.Range("I94").Value = 1670.2667
$ws.Range("J94").Value = 2172.2856
$ws.Range("K94").Value = 1670.2667
$ws.Range("L94").Value = 2172.2856
$ws.Range("M94").Value = -1219.2667
$ws.Range("N94").Value = -3074.2856

$ws.Range("H134").Value = 2974.3513
$ws.Range("I134").Value = 2854.5
$ws.Range("K134").Value = 8563.5
$ws.Range("M134").Value = -6028.5


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1535.8695
$ws.Range("J31").Value = 1711.5834
$ws.Range("L31").Value = 1711.5834
$ws.Range("N31").Value = -2301.5834

$ws.Range("H34").Value = 1535.8695
$ws.Range("J34").Value = 1711.5834
$ws.Range("L34").Value = 1711.5834
$ws.Range("N34").Value = -2115.5834

$ws.Range("H58").Value = 1938.5
$ws.Range("I58").Value = 1935.4
$ws.Range("K58").Value = 1935.4
$ws.Range("M58").Value = -1732.4

$ws.Range("H122").Value = 1970.8125
$ws.Range("I122").Value = 2058
$ws.Range("K122").Value = 6174
$ws.Range("M122").Value = -3724

$ws.Range("H132").Value = 2659.6765
$ws.Range("I132").Value = 2524.6897
$ws.Range("K132").Value = 7574.0691
$ws.Range("M132").Value = -5044.0691

$ws.Range("H134").Value = 2996.3333
$ws.Range("I134").Value = 2996.3333
$ws.Range("K134").Value = 8988.999899999999
$ws.Range("M134").Value = -6453.999899999999

$ws.Range("H136").Value = 1938.5
$ws.Range("I136").Value = 1935.4
$ws.Range("K136").Value = 5806.200000000001
$ws.Range("M136").Value = -3256.200000000001


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 4749.625
$ws.Range("J80").Value = 6249.5
$ws.Range("L80").Value = 18748.5
$ws.Range("N80").Value = -20620.5

$ws.Range("H83").Value = 4749.625
$ws.Range("J83").Value = 6249.5
$ws.Range("L83").Value = 56245.5
$ws.Range("N83").Value = -65605.5


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 21493.25
$ws.Range("I70").Value = 31329.6
$ws.Range("K70").Value = 31329.6
$ws.Range("M70").Value = -31059.6

$ws.Range("H73").Value = 21493.25
$ws.Range("I73").Value = 31329.6
$ws.Range("K73").Value = 31329.6
$ws.Range("M73").Value = -30393.6

$ws.Range("H113").Value = 2009.6316
$ws.Range("J113").Value = 2013.125
$ws.Range("L113").Value = 2013.125
$ws.Range("N113").Value = -6353.125

$ws.Range("H132").Value = 2361.5
$ws.Range("J132").Value = 2316.8
$ws.Range("L132").Value = 6950.400000000001
$ws.Range("N132").Value = -12010.4


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 974.7059
$ws.Range("I16").Value = 622.0714
$ws.Range("K16").Value = 622.0714
$ws.Range("M16").Value = -452.0714

$ws.Range("H22").Value = 1212.8605
$ws.Range("J22").Value = 1234.9487
$ws.Range("L22").Value = 1234.9487
$ws.Range("N22").Value = -1824.9487

$ws.Range("H27").Value = 1212.8605
$ws.Range("J27").Value = 1234.9487
$ws.Range("L27").Value = 1234.9487
$ws.Range("N27").Value = -1448.9487

$ws.Range("H122").Value = 9781.625
$ws.Range("I122").Value = 7791.3335
$ws.Range("K122").Value = 23374.0005
$ws.Range("M122").Value = -20924.0005

$ws.Range("H136").Value = 4680.654
$ws.Range("I136").Value = 4338.5654
$ws.Range("K136").Value = 13015.6962
$ws.Range("M136").Value = -10465.6962


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4227.5
$ws.Range("I81").Value = 1679.8
$ws.Range("K81").Value = 3359.6
$ws.Range("M81").Value = -2298.6

$ws.Range("H84").Value = 4227.5
$ws.Range("I84").Value = 1679.8
$ws.Range("K84").Value = 16798
$ws.Range("M84").Value = -11494

$ws.Range("H132").Value = 1082.1052
$ws.Range("J132").Value = 1161.25
$ws.Range("L132").Value = 3483.75
$ws.Range("N132").Value = -8543.75

$ws.Range("H136").Value = 3166.1177
$ws.Range("I136").Value = 2883.5
$ws.Range("J136").Value = 4485
$ws.Range("K136").Value = 8650.5
$ws.Range("L136").Value = 13455
$ws.Range("M136").Value = -6100.5
$ws.Range("N136").Value = -18555

